$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate formatting of an existing "Pass" row (row 4) into the new row 14,
# then overwrite with the new test case's data.
$ws.Range("A4:F4").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)

$ws.Range("A14").Value2 = "DeleteAccount"
$ws.Range("B14").Value2 = "A user should be able to delete their account and try to login once it has been deleted"
$ws.Range("C14").Value2 = 11
$ws.Range("D14").Value2 = "Pass"
$ws.Range("E14").Value2 = 42102
$ws.Range("F14").Value2 = 0.6

$ws.Range("E23").Select()
